$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column A (duplicate GENE values with border style), which
# shifts columns B:F left to A:E.
$ws.Range("A:A").EntireColumn.Delete()

# Fix the header text that was mis-spelled as MODEL_CONDITION; after the
# column shift it now lives in D1.
$ws.Range("D1").Value = "MODELCONDITION"
